# D14's shared string text itself changed ("... a = 0.01" -> "... a = 1/100");
# re-assert the value so the shared string entry is updated accordingly.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CI")

$ws.Range("D14").Value = "confidence of 99% impplies a = 1/100"

# E14: new content (previously empty)
$ws.Range("E14").Value = "a = 0.01"

# D15: updated label text ("a" -> "Z/2  = 0.01/2"); E15 value now 0.005 (was 0.01)
$ws.Range("D15").Value = "Z/2  = 0.01/2"
$ws.Range("E15").Value = 0.005

# D16: updated label text ("1-0.001 = " -> "1-0.005"); E16 value unchanged (0.995)
$ws.Range("D16").Value = "1-0.005"

# D17 / E17: new label + value (previously empty)
$ws.Range("D17").Value = "Z0.005 = 0.995"
$ws.Range("E17").Value = 0.05

# D19 / E19: updated label + t-score value (2.58 -> 2.57)
$ws.Range("D19").Value = "Za = 2.5+0.07"
$ws.Range("E19").Value = 2.57

# Update the saved selection to match the authored workbook state
$ws.Range("F19").Select()

$wb.Save()
